# "Changing waits in suite B"
# Update the Results column (E) for the "Test Cases" sheet (suite B rows)
# and refresh the active selection/view state on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# TestCase_B1 result: FAIL -> SKIP
$ws.Range("E2").Value = "SKIP"

# TestCase_B10 result: SKIP -> FAIL
$ws.Range("E11").Value = "FAIL"

# TestCase_B68 (row69) and TestCase_B83 (row84) results go from PASS -> SKIP,
# and TestCase_B84..TestCase_B88 (rows 85-89) were blank and now become SKIP.
$ws.Range("E69:E89").Value = "SKIP"

# Update the sheet's view/selection state (was topLeftCell A87 / selection C87)
$ws.Activate()
$ws.Range("D2:D89").Select()
